# Fixed a bug in removeSymbols
# This script re-orders the data rows (rows 2-25, columns A-F) of the active
# worksheet according to the permutation captured by the diff. Row 1 (header)
# and row 26 (totals) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (both are 1-based worksheet rows)
$mapping = @{
    2  = 3
    3  = 5
    4  = 14
    5  = 15
    6  = 9
    7  = 4
    8  = 6
    9  = 7
    10 = 2
    11 = 11
    12 = 8
    13 = 10
    14 = 13
    15 = 12
    16 = 20
    17 = 17
    18 = 21
    19 = 19
    20 = 16
    21 = 18
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

# Snapshot the original values for rows 2-25, columns A-F before any writes,
# so that overlapping writes don't clobber source data we still need to read.
$original = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $original[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $original[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
